$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the three shared strings that mention "UMR 1283" so they all
#    become "UMR 1283/8199" (commit: "add double number for institution").
#    Updating the cell .Value directly lets the engine repack the shared
#    string table itself (matches the observed before/after diff exactly).
# ---------------------------------------------------------------------------

# "experience" sheet (activities paragraph used in row 8 / B8)
$wsExperience = $wb.Worksheets.Item("experience")
$wsExperience.Range("B8").Value = "Functional (Epi)genomics and Molecular  `r`nPhysiology of Diabetes and Related Diseases  `r`nEGID - UMR 1283/8199  `r`n(European Genomics Institute for Diabetes)"

# "profil" sheet (paragraph describing current position, C2)
$wsProfil = $wb.Worksheets.Item("profil")
$wsProfil.Range("C2").Value = "I'm currently working for the *Institut Pasteur de Lille* at the *UMR 1283/8199 - ""Functional (Epi)genomics and Molecular Physiology of Diabetes and Related Diseases""* as the *head of the biostatistic team*."

# "contact" sheet (institute markdown link, B2)
$wsContact = $wb.Worksheets.Item("contact")
$wsContact.Range("B2").Value = "[EGID - UMR 1283/8199](http://www.good.cnrs.fr/?lang=en)"

# ---------------------------------------------------------------------------
# 2. Move the active tab / selection from "experience" to "contact".
# ---------------------------------------------------------------------------

# "profil" sheet keeps its own new selection (C2 -> C4) without changing tabs.
$wsProfil.Range("C4").Select() | Out-Null

# Switch to the "contact" sheet and select B12 there; this also flips which
# sheetView carries tabSelected="1" (off experience, on contact) and updates
# the workbook-level active tab.
$wsContact.Activate() | Out-Null
$wsContact.Range("B12").Select() | Out-Null

# "experience" sheet keeps its own selection at B8 (unchanged), it just loses
# the tabSelected flag because contact is now active (handled by Activate()).
$wsExperience.Range("B8").Select() | Out-Null

# Re-activate "contact" last so it is the sheet left active/selected.
$wsContact.Activate() | Out-Null
